# Refresh the cryptos worksheet with the latest scraped coinranking.com values
# (price / 1h volume-change snapshot, matching the GitHub Actions commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking prices such as "0.9991" or "305.30" must stay text (not become numbers),
# so force those specific Price cells to Text format before writing the value.
$priceTextCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D10",
    "D13",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D21",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.932.86'
$ws.Range("E2").Value = '  -2.79%  '
$ws.Range("D3").Value = '1.859.56'
$ws.Range("E3").Value = '  -2.20%  '
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '305.30'
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.5029'
$ws.Range("E7").Value = '  -2.83%  '
$ws.Range("E8").Value = '  -1.67%  '
$ws.Range("E9").Value = '  -1.68%  '
$ws.Range("D10").Value = '0.8849'
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("E11").Value = '  -2.83%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.876.25'
$ws.Range("E12").Value = '  -1.34%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07558'
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("D14").Value = '5.280'
$ws.Range("E14").Value = '  -2.98%  '
$ws.Range("D15").Value = '88.97'
$ws.Range("E15").Value = '  -3.27%  '
$ws.Range("D16").Value = '0.9993'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '0.000008347'
$ws.Range("E17").Value = '  -4.26%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Value = '14.04'
$ws.Range("E18").Value = '  -2.89%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '0.9988'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '26.983.51'
$ws.Range("E20").Value = '  -2.67%  '
$ws.Range("D21").Value = '5.021'
$ws.Range("E21").Value = '  -2.18%  '
$ws.Range("D22").Value = '2.113.28'
$ws.Range("E22").Value = '  -1.30%  '
$ws.Range("D23").Value = '10.44'
$ws.Range("E23").Value = '  -3.48%  '
$ws.Range("D24").Value = '6.451'
$ws.Range("E24").Value = '  -1.97%  '
$ws.Range("D25").Value = '1.847'
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("D26").Value = '146.68'
$ws.Range("E26").Value = '  -4.43%  '
$ws.Range("D27").Value = '17.90'
$ws.Range("E27").Value = '  -2.24%  '
$ws.Range("D28").Value = '2.090'
$ws.Range("E28").Value = '  -4.44%  '
$ws.Range("D29").Value = '112.29'
$ws.Range("E29").Value = '  -2.23%  '
$ws.Range("D30").Value = '4.634'
$ws.Range("E30").Value = '  -4.40%  '
$ws.Range("D31").Value = '4.649'
$ws.Range("E31").Value = '  -3.54%  '
$ws.Range("D32").Value = '0.09029'
$ws.Range("E32").Value = '  +0.51%  '
$ws.Range("D33").Value = '0.05100'
$ws.Range("E33").Value = '  -2.90%  '
$ws.Range("D34").Value = '3.040'
$ws.Range("E34").Value = '  -4.27%  '
$ws.Range("D35").Value = '1.148'
$ws.Range("E35").Value = '  -7.24%  '
$ws.Range("D36").Value = '0.7213'
$ws.Range("E36").Value = '  -7.52%  '
$ws.Range("D37").Value = '0.02030'
$ws.Range("E37").Value = '  -2.44%  '
$ws.Range("D38").Value = '3.035'
$ws.Range("E38").Value = '  -0.56%  '
$ws.Range("D39").Value = '2.453'
$ws.Range("E39").Value = '  -6.17%  '
$ws.Range("D40").Value = '1.071'
$ws.Range("E40").Value = '  -1.65%  '
$ws.Range("D41").Value = '0.5268'
$ws.Range("E41").Value = '  -4.16%  '
$ws.Range("D42").Value = '6.521'
$ws.Range("E42").Value = '  -2.20%  '
$ws.Range("D43").Value = '115.02'
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("D44").Value = '8.226'
$ws.Range("E44").Value = '  -3.13%  '
$ws.Range("D45").Value = '0.1460'
$ws.Range("E45").Value = '  -3.02%  '
$ws.Range("D46").Value = '0.9986'
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = '0.4582'
$ws.Range("E47").Value = '  -4.43%  '
$ws.Range("D48").Value = '9.948'
$ws.Range("E48").Value = '  -4.32%  '
$ws.Range("D49").Value = '1.554'
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("D50").Value = '36.37'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").Value = '63.82'
$ws.Range("E51").Value = '  -4.13%  '
